$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 25 and 26) since the table now has 23 data rows instead of 25
$ws.Rows.Item(25).EntireRow.Delete() | Out-Null
$ws.Rows.Item(25).EntireRow.Delete() | Out-Null

# Update remaining rows (2-24) to reflect the refreshed GoodInfo stock list data
# Row 2
$ws.Range("A2").Value = 64

# Row 3
$ws.Range("A3").Value = 65
$ws.Range("C3").Value = 2484
$ws.Range("D3").Value = "希華"
$ws.Range("F3").Value = 11
$ws.Range("H3").Value = 1.86
$ws.Range("J3").Value = 3.84
$ws.Range("K3").Value = 1.86
$ws.Range("L3").Value = 0.4843750000000001

# Row 4
$ws.Range("A4").Value = 66
$ws.Range("C4").Value = 5351
$ws.Range("D4").Value = "鈺創"
$ws.Range("F4").Value = 37
$ws.Range("H4").Value = 0.7
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 0.7
$ws.Range("L4").Value = 0.6363636363636362

# Row 5
$ws.Range("A5").Value = 67
$ws.Range("C5").Value = 6265
$ws.Range("D5").Value = "方土昶"
$ws.Range("F5").Value = 49
$ws.Range("H5").Value = 0.75
$ws.Range("J5").Value = 2.76
$ws.Range("K5").Value = 0.75
$ws.Range("L5").Value = 0.2717391304347826

# Row 6
$ws.Range("A6").Value = 68
$ws.Range("C6").Value = 6568
$ws.Range("D6").Value = "宏觀"
$ws.Range("F6").Value = 54
$ws.Range("H6").Value = 2.12
$ws.Range("I6").Value = 11.7
$ws.Range("J6").Value = $null
$ws.Range("K6").Value = 13.82
$ws.Range("L6").Value = $null

# Row 7
$ws.Range("A7").Value = 69
$ws.Range("C7").Value = 3041
$ws.Range("D7").Value = "揚智"
$ws.Range("F7").Value = 17
$ws.Range("G7").Value = 2.9
$ws.Range("H7").Value = 22.8
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = 16.25
$ws.Range("K7").Value = 22.8
$ws.Range("L7").Value = 1.403076923076923

# Row 8
$ws.Range("A8").Value = 70
$ws.Range("C8").Value = 3122
$ws.Range("D8").Value = "笙泉"
$ws.Range("F8").Value = 20
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 9.14
$ws.Range("J8").Value = 12.34
$ws.Range("K8").Value = 12.34
$ws.Range("L8").Value = 1

# Row 9
$ws.Range("A9").Value = 71
$ws.Range("C9").Value = 3221
$ws.Range("D9").Value = "台嘉碩"
$ws.Range("F9").Value = 24
$ws.Range("H9").Value = 7.84
$ws.Range("I9").Value = 3.8
$ws.Range("J9").Value = $null
$ws.Range("K9").Value = 11.64
$ws.Range("L9").Value = $null

# Row 10
$ws.Range("A10").Value = 72
$ws.Range("C10").Value = 3515
$ws.Range("D10").Value = "華擎"
$ws.Range("F10").Value = 30
$ws.Range("H10").Value = 4.24
$ws.Range("I10").Value = 11.24
$ws.Range("J10").Value = 4.17
$ws.Range("K10").Value = 15.48
$ws.Range("L10").Value = 3.712230215827338

# Row 11
$ws.Range("A11").Value = 73
$ws.Range("C11").Value = 3588
$ws.Range("D11").Value = "通嘉"
$ws.Range("F11").Value = 32
$ws.Range("H11").Value = 3.22
$ws.Range("I11").Value = 5.96
$ws.Range("J11").Value = 8.18
$ws.Range("K11").Value = 9.18
$ws.Range("L11").Value = 1.122249388753056

# Row 12
$ws.Range("A12").Value = 74
$ws.Range("C12").Value = 5371
$ws.Range("D12").Value = "中光電"
$ws.Range("F12").Value = 38
$ws.Range("H12").Value = 17.58
$ws.Range("I12").Value = $null
$ws.Range("J12").Value = 11.37
$ws.Range("K12").Value = 17.58
$ws.Range("L12").Value = 1.546174142480211

# Row 13
$ws.Range("A13").Value = 75
$ws.Range("C13").Value = 6411
$ws.Range("D13").Value = "晶焱"
$ws.Range("F13").Value = 52
$ws.Range("H13").Value = 0.56
$ws.Range("I13").Value = 5.92
$ws.Range("J13").Value = $null
$ws.Range("K13").Value = 6.48
$ws.Range("L13").Value = $null

# Row 14
$ws.Range("A14").Value = 76
$ws.Range("C14").Value = 8040
$ws.Range("D14").Value = "九暘"
$ws.Range("F14").Value = 57
$ws.Range("H14").Value = 7.96
$ws.Range("I14").Value = 9.2
$ws.Range("K14").Value = 17.16

# Row 15
$ws.Range("A15").Value = 77
$ws.Range("C15").Value = 1110
$ws.Range("D15").Value = "東泥"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2.5
$ws.Range("H15").Value = 4.46
$ws.Range("I15").Value = $null
$ws.Range("K15").Value = 4.46

# Row 16
$ws.Range("A16").Value = 78
$ws.Range("C16").Value = 2338
$ws.Range("D16").Value = "光罩"
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = $null
$ws.Range("K16").Value = 0

# Row 17
$ws.Range("A17").Value = 79
$ws.Range("C17").Value = 3312
$ws.Range("D17").Value = "弘憶股"
$ws.Range("F17").Value = 26

# Row 18
$ws.Range("A18").Value = 80
$ws.Range("C18").Value = 3455
$ws.Range("D18").Value = "由田"
$ws.Range("F18").Value = 28

# Row 19
$ws.Range("A19").Value = 81
$ws.Range("C19").Value = 5258
$ws.Range("D19").Value = "虹堡"
$ws.Range("F19").Value = 35

# Row 20
$ws.Range("A20").Value = 82
$ws.Range("C20").Value = 5347
$ws.Range("D20").Value = "世界"
$ws.Range("F20").Value = 36

# Row 21
$ws.Range("A21").Value = 83
$ws.Range("C21").Value = 6143
$ws.Range("D21").Value = "振曜"
$ws.Range("F21").Value = 43

# Row 22
$ws.Range("A22").Value = 84
$ws.Range("C22").Value = 6285
$ws.Range("D22").Value = "啟碁"
$ws.Range("F22").Value = 51

# Row 23
$ws.Range("A23").Value = 85
$ws.Range("C23").Value = 6438
$ws.Range("D23").Value = "迅得"
$ws.Range("F23").Value = 53

# Row 24
$ws.Range("A24").Value = 86
$ws.Range("C24").Value = 8054
$ws.Range("D24").Value = "安國"
$ws.Range("F24").Value = 59
